$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.746.03'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').Value = '''1.888.81'
$ws.Range('E3').Value = '  -0.85%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '''0.7936'
$ws.Range('E5').Value = '  -1.57%  '
$ws.Range('D6').Value = '''241.61'
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '''0.3171'
$ws.Range('E8').Value = '  +1.62%  '
$ws.Range('D9').Value = '''25.57'
$ws.Range('E9').Value = '  -3.40%  '
$ws.Range('D10').Value = '''0.07033'
$ws.Range('E10').Value = '  +0.09%  '
$ws.Range('D11').Value = '''0.08048'
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('D12').Value = '''0.7673'
$ws.Range('E12').Value = '  +3.28%  '
$ws.Range('D13').Value = '''1.908.52'
$ws.Range('E13').Value = '  +0.19%  '
$ws.Range('D14').Value = '''5.316'
$ws.Range('E14').Value = '  +2.62%  '
$ws.Range('D15').Value = '''92.02'
$ws.Range('E15').Value = '  -0.42%  '
$ws.Range('D16').Value = '''29.752.88'
$ws.Range('E16').Value = '  -0.66%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').Value = '''13.78'
$ws.Range('E17').Value = '  -1.45%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = '''5.930'
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('D19').Value = '''242.91'
$ws.Range('E19').Value = '  -1.00%  '
$ws.Range('D20').Value = '''0.000007686'
$ws.Range('E20').Value = '  -1.19%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').Value = '''8.219'
$ws.Range('E21').Value = '  +18.64%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '''1.000'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').Value = '''2.138.06'
$ws.Range('E23').Value = '  -0.88%  '
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('D25').Value = '''0.1625'
$ws.Range('E25').Value = '  +9.63%  '
$ws.Range('D26').Value = '''9.303'
$ws.Range('E26').Value = '  +1.07%  '
$ws.Range('D27').Value = '''163.73'
$ws.Range('E27').Value = '  -2.82%  '
$ws.Range('D28').Value = '''18.65'
$ws.Range('E28').Value = '  -1.21%  '
$ws.Range('D29').Value = '''2.051'
$ws.Range('E29').Value = '  -0.78%  '
$ws.Range('D30').Value = '''1.374'
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('D31').Value = '''1.533'
$ws.Range('E31').Value = '  +1.47%  '
$ws.Range('D32').Value = '''4.437'
$ws.Range('E32').Value = '  +3.41%  '
$ws.Range('D33').Value = '''0.05692'
$ws.Range('E33').Value = '  +2.86%  '
$ws.Range('D34').Value = '''4.083'
$ws.Range('E34').Value = '  +0.61%  '
$ws.Range('D35').Value = '''1.262'
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('D36').Value = '''0.7369'
$ws.Range('E36').Value = '  +1.08%  '
$ws.Range('D37').Value = '''1.001'
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('D38').Value = '''2.705'
$ws.Range('E38').Value = '  -0.29%  '
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('D40').Value = '''2.770'
$ws.Range('E40').Value = '  -0.49%  '
$ws.Range('D41').Value = '''0.4401'
$ws.Range('E41').Value = '  -0.24%  '
$ws.Range('D42').Value = '''72.30'
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('D43').Value = '''5.833'
$ws.Range('E43').Value = '  -2.35%  '
$ws.Range('D44').Value = '''1.000'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').Value = '''0.8394'
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('D46').Value = '''1.030.50'
$ws.Range('E46').Value = '  +4.52%  '
$ws.Range('D47').Value = '''102.08'
$ws.Range('E47').Value = '  +1.22%  '
$ws.Range('D48').Value = '''1.858'
$ws.Range('E48').Value = '  -1.57%  '
$ws.Range('D49').Value = '''9.863'
$ws.Range('E49').Value = '  +1.66%  '
$ws.Range('D50').Value = '''7.450'
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').Value = '''2.043.09'
$ws.Range('E51').Value = '  -0.75%  '
